# Generate Report for Handoff
# Update status + handoff datetime for the "75c8985e-..." row (row 3) across
# the Overview, zh-cn and de-de sheets to reflect the file now being ready
# for handoff again.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: B3 = zh-cn status, C3 = de-de status, D3 = Latest Handoff Date
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-40-12 20:40:20"

# zh-cn sheet: C3 = Status, E3 = Latest Handoff Datetime
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-12 20:40:16"

# de-de sheet: C3 = Status, E3 = Latest Handoff Datetime
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-12 20:40:20"
